$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30..79 down to 31..80
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value2 = 45174
$ws.Range("D30").NumberFormat = $ws.Range("D31").NumberFormat
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100104
$ws.Range("H30").Value = "Frutos de pepita"
$ws.Range("I30").Value = 100104005
$ws.Range("J30").Value = "Pera"
$ws.Range("K30").Value = "Packham's Triumph"
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 21000
$ws.Range("O30").Value = 22000
$ws.Range("P30").Value = 21500
$ws.Range("Q30").Value = "$/caja 18 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 1194
$ws.Range("T30").Value = 18
